$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A9 previously held the "last row" date format (YYYY-MM-DD); since it is no
# longer the last row, give it the same format as the other data rows
# (YYYY-MM-DD HH:MM:SS), matching A2:A8.
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new data row for 2021-11-19.
$ws.Range("A10").Value = 44519
$ws.Range("B10").Value = -138.0499999999997

# The newly appended row becomes the new "last row", so it gets the
# YYYY-MM-DD only format (same as A9 had before the edit).
$ws.Range("A10").NumberFormat = "YYYY-MM-DD"
